$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Data edits ---
# Row 29: "Integrar API para cadastro e login" item size corrected from M(13) to G(21)
$ws.Range("E29").Value = 21

# N5 / N6: fill in previously empty WEB-DATA-VIZ sprint totals
$ws.Range("N5").Value = 13
$ws.Range("N6").Value = 21

# --- Highlight rows 22, 23, 29 (mark as integrated / moved to WEB-DATA-VIZ) ---
# Copy the fill/format already used by completed rows (e.g. row 9) onto these rows.
$ws.Range("A9:H9").Copy()
$ws.Range("A22:H22").PasteSpecial(-4122)
$ws.Range("A23:H23").PasteSpecial(-4122)
$ws.Range("A29:H29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Cursor / selection change ---
$ws.Range("P4").Select()
